# Adding new FDOM and POC data
# Appends 8 new sample rows (80-87) to Sheet1, matching the pattern of the
# existing data table (columns A-K), extends the two shared formulas in
# columns J and K down through the new rows, and copies the date-number
# formatting from the existing rows so the new cells use the same styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new data rows 80-87 (FDOM and POC data) ---

# Row 80
$ws.Cells.Item(80,1).Value = 45505
$ws.Cells.Item(80,2).Value = 6.6
$ws.Cells.Item(80,3).Value = 565
$ws.Cells.Item(80,4).Value = 1.0686
$ws.Cells.Item(80,5).Value = 45509
$ws.Cells.Item(80,6).Value = 1.1885
$ws.Cells.Item(80,7).Value = 45509
$ws.Cells.Item(80,8).Value = 1.1814
$ws.Cells.Item(80,9).Value = "AB"

# Row 81
$ws.Cells.Item(81,1).Value = 45502
$ws.Cells.Item(81,2).Value = 15
$ws.Cells.Item(81,3).Value = 560
$ws.Cells.Item(81,4).Value = 1.1224
$ws.Cells.Item(81,5).Value = 45509
$ws.Cells.Item(81,6).Value = 1.2424
$ws.Cells.Item(81,7).Value = 45509
$ws.Cells.Item(81,8).Value = 1.2336
$ws.Cells.Item(81,9).Value = "AB"

# Row 82
$ws.Cells.Item(82,1).Value = 45502
$ws.Cells.Item(82,2).Value = 9
$ws.Cells.Item(82,3).Value = 563
$ws.Cells.Item(82,4).Value = 1.1046
$ws.Cells.Item(82,5).Value = 45509
$ws.Cells.Item(82,6).Value = 1.2242
$ws.Cells.Item(82,7).Value = 45509
$ws.Cells.Item(82,8).Value = 1.2197
$ws.Cells.Item(82,9).Value = "AB"

# Row 83
$ws.Cells.Item(83,1).Value = 45502
$ws.Cells.Item(83,2).Value = 6
$ws.Cells.Item(83,3).Value = 554
$ws.Cells.Item(83,4).Value = 1.1166
$ws.Cells.Item(83,5).Value = 45509
$ws.Cells.Item(83,6).Value = 1.2376
$ws.Cells.Item(83,7).Value = 45509
$ws.Cells.Item(83,8).Value = 1.2316
$ws.Cells.Item(83,9).Value = "AB"

# Row 84
$ws.Cells.Item(84,1).Value = 45505
$ws.Cells.Item(84,2).Value = 5.5
$ws.Cells.Item(84,3).Value = 556
$ws.Cells.Item(84,4).Value = 1.0638
$ws.Cells.Item(84,5).Value = 45509
$ws.Cells.Item(84,6).Value = 1.1865
$ws.Cells.Item(84,7).Value = 45509
$ws.Cells.Item(84,8).Value = 1.18
$ws.Cells.Item(84,9).Value = "AB"

# Row 85 (ID column holds text "6a" instead of a number)
$ws.Cells.Item(85,1).Value = 45502
$ws.Cells.Item(85,2).Value = "6a"
$ws.Cells.Item(85,3).Value = 563
$ws.Cells.Item(85,4).Value = 1.1359
$ws.Cells.Item(85,5).Value = 45509
$ws.Cells.Item(85,6).Value = 1.2546
$ws.Cells.Item(85,7).Value = 45509
$ws.Cells.Item(85,8).Value = 1.2508
$ws.Cells.Item(85,9).Value = "AB"

# Row 86
$ws.Cells.Item(86,1).Value = 45502
$ws.Cells.Item(86,2).Value = 5
$ws.Cells.Item(86,3).Value = 564
$ws.Cells.Item(86,4).Value = 1.1136
$ws.Cells.Item(86,5).Value = 45509
$ws.Cells.Item(86,6).Value = 1.2343
$ws.Cells.Item(86,7).Value = 45509
$ws.Cells.Item(86,8).Value = 1.2286
$ws.Cells.Item(86,9).Value = "AB"

# Row 87
$ws.Cells.Item(87,1).Value = 45505
$ws.Cells.Item(87,2).Value = 9.5
$ws.Cells.Item(87,3).Value = 562
$ws.Cells.Item(87,4).Value = 1.1169
$ws.Cells.Item(87,5).Value = 45509
$ws.Cells.Item(87,6).Value = 1.2407
$ws.Cells.Item(87,7).Value = 45509
$ws.Cells.Item(87,8).Value = 1.2341
$ws.Cells.Item(87,9).Value = "AB"

# --- Re-use the existing date-number-format styles (columns A, E, G use the
#     same "Sampled"/"Desicator"/"Furnace" date styles as the rest of the
#     table) by copying formats only, so no new style/numFmt entries are
#     created in styles.xml ---
$ws.Range("A2").Copy()
$ws.Range("A80:A87").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("E2").Copy()
$ws.Range("E80:E87").PasteSpecial(-4122)

$ws.Range("E2").Copy()
$ws.Range("G80:G87").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Extend the shared formulas in columns J and K down through row 87 ---
$ws.Range("J4:J87").Formula = "=(F4-H4)*1000"
$ws.Range("K2:K87").Formula = "=J2/(C2/1000)"

# --- Restore the selection to match the new working location ---
$ws.Range("O90").Select()
